$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.15
$summary.Range("B6").Value = 28
$summary.Range("B9").Value = 42.86

# --- Sheet: Strategy Status ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 28
$status.Range("G5").Value = 42.86

# --- New trade row data (Trade #28), appended to "All Trades" and "MarketMaking" ---
function Add-TradeRow($ws, $rowNum) {
    $ws.Cells.Item($rowNum, 1).Value = 28
    # The date-like string needs to be forced as text, otherwise Excel
    # auto-converts "2026-02-17" into a date serial number.
    $dateCell = $ws.Cells.Item($rowNum, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()
    $ws.Cells.Item($rowNum, 3).Value = "20:07:59"
    $ws.Cells.Item($rowNum, 4).Value = "MarketMaking"
    $ws.Cells.Item($rowNum, 5).Value = "DOWN"
    $ws.Cells.Item($rowNum, 6).Value = 0.01
    $ws.Cells.Item($rowNum, 7).Value = 0.01
    $ws.Cells.Item($rowNum, 8).Value = "CLOSED"
    $ws.Cells.Item($rowNum, 9).Value = 0
    $ws.Cells.Item($rowNum, 10).Value = 0
    $ws.Cells.Item($rowNum, 11).Value = 99.8
    $ws.Cells.Item($rowNum, 12).Value = 0
    $ws.Cells.Item($rowNum, 13).Value = 0
    $ws.Cells.Item($rowNum, 14).Value = 0.6
    $ws.Cells.Item($rowNum, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($rowNum, 16).Value = "early_exit"
    $ws.Cells.Item($rowNum, 17).Value = 0.11
}

# --- Sheet: All Trades ---
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 29

# --- Sheet: MarketMaking ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 29
